{"js": "// Apply text replacements for the multiplication table answers.\n// Each pair is [oldText, newText] taken from the diff.\nconst pairs = [\n  [\"69\u00d726=1794\", \"94\u00d724=2256\"],\n  [\"89\u00d733=2937\", \"98\u00d792=9016\"],\n  [\"47\u00d727=1269\", \"91\u00d786=7826\"],\n  [\"38\u00d735=1330\", \"13\u00d763=819\"],\n  [\"60\u00d769=4140\", \"30\u00d755=1650\"],\n  [\"77\u00d719=1463\", \"38\u00d746=1748\"],\n  [\"33\u00d787=2871\", \"87\u00d728=2436\"],\n  [\"72\u00d728=2016\", \"86\u00d759=5074\"],\n  [\"40\u00d771=2840\", \"51\u00d733=1683\"],\n  [\"64\u00d799=6336\", \"99\u00d791=9009\"],\n  [\"17\u00d775=1275\", \"52\u00d773=3796\"],\n  [\"78\u00d744=3432\", \"27\u00d745=1215\"],\n  [\"64\u00d781=5184\", \"91\u00d791=8281\"],\n  [\"93\u00d727=2511\", \"49\u00d741=2009\"],\n  [\"74\u00d780=5920\", \"58\u00d764=3712\"],\n  [\"84\u00d724=2016\", \"96\u00d799=9504\"],\n  [\"87\u00d722=1914\", \"51\u00d777=3927\"],\n  [\"52\u00d771=3692\", \"43\u00d772=3096\"],\n  [\"34\u00d729=986\", \"44\u00d786=3784\"],\n  [\"85\u00d790=7650\", \"38\u00d779=3002\"],\n  [\"94\u00d748=4512\", \"23\u00d760=1380\"],\n  [\"72\u00d714=1008\", \"82\u00d774=6068\"],\n  [\"70\u00d751=3570\", \"39\u00d760=2340\"],\n  [\"65\u00d781=5265\", \"76\u00d723=1748\"],\n  [\"38\u00d751=1938\", \"85\u00d770=5950\"]\n];\n\n// Step 1: search for every old value in the document body.\nconst searchResults = [];\nfor (const [oldText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\n// Step 2: replace each found range with its new value.\nfor (let i = 0; i < pairs.length; i++) {\n  const newText = pairs[i][1];\n  const results = searchResults[i];\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply text replacements for the multiplication table answers.\n# Each pair is (oldText, newText) taken from the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"69\u00d726=1794\", \"94\u00d724=2256\"),\n    @(\"89\u00d733=2937\", \"98\u00d792=9016\"),\n    @(\"47\u00d727=1269\", \"91\u00d786=7826\"),\n    @(\"38\u00d735=1330\", \"13\u00d763=819\"),\n    @(\"60\u00d769=4140\", \"30\u00d755=1650\"),\n    @(\"77\u00d719=1463\", \"38\u00d746=1748\"),\n    @(\"33\u00d787=2871\", \"87\u00d728=2436\"),\n    @(\"72\u00d728=2016\", \"86\u00d759=5074\"),\n    @(\"40\u00d771=2840\", \"51\u00d733=1683\"),\n    @(\"64\u00d799=6336\", \"99\u00d791=9009\"),\n    @(\"17\u00d775=1275\", \"52\u00d773=3796\"),\n    @(\"78\u00d744=3432\", \"27\u00d745=1215\"),\n    @(\"64\u00d781=5184\", \"91\u00d791=8281\"),\n    @(\"93\u00d727=2511\", \"49\u00d741=2009\"),\n    @(\"74\u00d780=5920\", \"58\u00d764=3712\"),\n    @(\"84\u00d724=2016\", \"96\u00d799=9504\"),\n    @(\"87\u00d722=1914\", \"51\u00d777=3927\"),\n    @(\"52\u00d771=3692\", \"43\u00d772=3096\"),\n    @(\"34\u00d729=986\", \"44\u00d786=3784\"),\n    @(\"85\u00d790=7650\", \"38\u00d779=3002\"),\n    @(\"94\u00d748=4512\", \"23\u00d760=1380\"),\n    @(\"72\u00d714=1008\", \"82\u00d774=6068\"),\n    @(\"70\u00d751=3570\", \"39\u00d760=2340\"),\n    @(\"65\u00d781=5265\", \"76\u00d723=1748\"),\n    @(\"38\u00d751=1938\", \"85\u00d770=5950\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found: $oldText\"\n    }\n}\n"}
